$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the existing "Z_xsec" column header (col L) to "Z_xsec_norm"
$ws.Range("L2:L6").Value = "Z_xsec_norm"

# Add new column M: "normalization" header + data values (7.7 for each data row)
$ws.Range("M1").Value = "normalization"
$ws.Range("M2").Value = 7.7
$ws.Range("M3").Value = 7.7
$ws.Range("M4").Value = 7.7
$ws.Range("M5").Value = 7.7
$ws.Range("M6").Value = 7.7

# New column M cells share the same centered style as the rest of the table
$ws.Range("M1:M6").HorizontalAlignment = -4108

# Widen the new column to fit its content
$ws.Columns("M").ColumnWidth = 12.8

# Move the active selection to reflect where editing finished
$ws.Range("M7").Select()
